# Generate Report for Handback
# Applies:
#  - Overview sheet: status text "In Translation" -> "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: populate "Latest Target File" (I) and "Latest Handback File" (J)
#    columns (and "Latest Handback DateTime" (K) for de-de) with hyperlinked file names,
#    matching the style already used for the "Source File Name" (A) column.
#  - Widen a few columns that now hold long file names / urls.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$cornflowerBlue = 15570276   # BGR value for RGB FF6495ED (matches existing HyperLink font color)

function Set-Handback($ws, $iCell, $jCell, $fileName, $xlfName, $url) {
    # "Latest Target File" column: source file name, hyperlinked like column A
    $ws.Range($iCell).Value = $fileName
    $ws.Hyperlinks.Add($ws.Range($iCell), $url, "", "", $fileName)
    $ws.Range($iCell).Font.Underline = 2
    $ws.Range($iCell).Font.Color = $cornflowerBlue
    $ws.Range($iCell).Font.Name = "Calibri"
    $ws.Range($iCell).Font.Size = 11

    # "Latest Handback File" column: generated xlf file name (plain text)
    $ws.Range($jCell).Value = $xlfName
}

# ---------------------------------------------------------------------------
# Overview sheet: update status cells
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.14
$overview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Columns.Item(3).ColumnWidth = 29.14
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

Set-Handback $zhcn "I2" "J2" `
    "3b1a8d4a-dba2-47e6-934a-23507a94d9a3.md" `
    "3b1a8d4a-dba2-47e6-934a-23507a94d9a3.c86e789a20c9b687dd8f9ad3814739373963960b.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c88aa4c14a52e4573d063da2b0c6156f473690/e2e/3b1a8d4a-dba2-47e6-934a-23507a94d9a3.md"

Set-Handback $zhcn "I3" "J3" `
    "6f54b17e-5c03-4f91-a3fb-ea82d0a0221c.yml" `
    "6f54b17e-5c03-4f91-a3fb-ea82d0a0221c.077971ed7a3e39da4e47dc9ca4a4224a6d48c0f7.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c88aa4c14a52e4573d063da2b0c6156f473690/e2e/6f54b17e-5c03-4f91-a3fb-ea82d0a0221c.yml"

Set-Handback $zhcn "I4" "J4" `
    "e650ad6b-e602-481e-9e95-1e607b783c7b.yml" `
    "e650ad6b-e602-481e-9e95-1e607b783c7b.4232caac424a820f88d29a25e4d13ae2a57905bb.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c88aa4c14a52e4573d063da2b0c6156f473690/e2e/e650ad6b-e602-481e-9e95-1e607b783c7b.yml"

# "Latest Handback DateTime" for the zh-cn rows picks up the globally updated
# "0001-01-01 00:00:00" -> "2016-11-14 08:00:23" text.
$zhcn.Range("K2").Value = "2016-11-14 08:00:23"
$zhcn.Range("K3").Value = "2016-11-14 08:00:23"
$zhcn.Range("K4").Value = "2016-11-14 08:00:23"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Columns.Item(3).ColumnWidth = 29.14
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17

Set-Handback $dede "I2" "J2" `
    "3b1a8d4a-dba2-47e6-934a-23507a94d9a3.md" `
    "3b1a8d4a-dba2-47e6-934a-23507a94d9a3.c86e789a20c9b687dd8f9ad3814739373963960b.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c88aa4c14a52e4573d063da2b0c6156f473690/e2e/3b1a8d4a-dba2-47e6-934a-23507a94d9a3.md"
$dede.Range("K2").Value = "2016-11-14 08:00:42"

Set-Handback $dede "I3" "J3" `
    "6f54b17e-5c03-4f91-a3fb-ea82d0a0221c.yml" `
    "6f54b17e-5c03-4f91-a3fb-ea82d0a0221c.077971ed7a3e39da4e47dc9ca4a4224a6d48c0f7.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c88aa4c14a52e4573d063da2b0c6156f473690/e2e/6f54b17e-5c03-4f91-a3fb-ea82d0a0221c.yml"
$dede.Range("K3").Value = "2016-11-14 08:00:42"

Set-Handback $dede "I4" "J4" `
    "e650ad6b-e602-481e-9e95-1e607b783c7b.yml" `
    "e650ad6b-e602-481e-9e95-1e607b783c7b.4232caac424a820f88d29a25e4d13ae2a57905bb.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c88aa4c14a52e4573d063da2b0c6156f473690/e2e/e650ad6b-e602-481e-9e95-1e607b783c7b.yml"
$dede.Range("K4").Value = "2016-11-14 08:00:42"
